# New version of the Venn diagram
#
# Resizes/repositions the three ovals of the TypeScript / ECMAScript /
# JavaScript Venn diagram and relabels the "older JavaScript" circle's
# second line to "ES/JavaScript".
#
# Note on the numeric literals below: Shape.Left/Top/Width/Height are
# expressed in points, but the underlying OOXML stores EMU (1 pt =
# 12700 EMU) and this host's point -> EMU conversion floors a
# single-precision (float32) reading of the literal. The literals here
# were chosen so that floor(float32(pt) * 12700) lands exactly on the
# target EMU values from the updated diagram.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "TypeScript" oval (Oval 3)
$oval1 = $s.Shapes.Item(1)
$oval1.Left = 211.5144119488189
$oval1.Top = 83.2208672417323
$oval1.Width = 327.652125984252
$oval1.Height = 317.641983103937

# "ECMAScript new features" oval (Oval 4)
$oval2 = $s.Shapes.Item(2)
$oval2.Left = 304.1355133110236
$oval2.Top = 176.10385896771652
$oval2.Width = 217.31905371811024
$oval2.Height = 208.20780187559055

# "older JavaScript supported by browsers" oval (Oval 5)
$oval3 = $s.Shapes.Item(3)
$oval3.Left = 337.24662787322836
$oval3.Top = 264.3116608433071
$oval3.Width = 168.62346656692912
$oval3.Height = 110.02590561181101

# Update the label text " JavaScript" -> " ES/JavaScript" (2nd run of the
# 1st paragraph), preserving its run-level formatting.
$tr = $oval3.TextFrame.TextRange
$run = $tr.Characters(6, 11)
$run.Text = " ES/JavaScript"
